# issue #5: stock data output to json file
#
# The "股票" (stock) sheet gains a new "property_category" column (always
# "stock") inserted right after the "total" column. The columns that used
# to hold the transaction date / legislator name / legislator id each shift
# one slot to the right, and a brand-new "legislator_id" column is appended
# at the far right (K). Two numeric-looking text cells that used to contain
# full-width commas ("2，286" / "4，797，260") are cleaned up to plain digit
# strings ("2286" / "4797260").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

$lastRow = 19

# Make sure the shifted "date" values (and the two cleaned-up numeric
# strings) keep being stored as text instead of Excel re-parsing them as
# dates/numbers when we write them back with .Value.
$ws.Range("I2:I" + $lastRow).NumberFormat = "@"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"

# --- header row -----------------------------------------------------------
# H1/I1/J1 keep their old shared-string slots (date/legislator_name/
# legislator_id in the pre-edit table); inserting "property_category" right
# before "date" in the string table makes H1 read "property_category" and
# pushes the rest down by one column. K1 is a brand new header cell.
$ws.Range("K1").Value = $ws.Range("J1").Value()
$ws.Range("K1").Font.Bold = $true
$ws.Range("K1").Borders.LineStyle = 1
$ws.Range("K1").HorizontalAlignment = -4108
$ws.Range("K1").VerticalAlignment = -4160
$ws.Range("J1").Value = $ws.Range("I1").Value()
$ws.Range("I1").Value = $ws.Range("H1").Value()
$ws.Range("H1").Value = "property_category"

# --- data rows --------------------------------------------------------------
for ($r = 2; $r -le $lastRow; $r++) {
    $oldH = $ws.Cells.Item($r, 8).Value()   # H: date, e.g. "2013-12-20"
    $oldI = $ws.Cells.Item($r, 9).Value()   # I: legislator name, e.g. "廖正井"
    $oldJ = $ws.Cells.Item($r, 10).Value()  # J: legislator id, e.g. 1711

    $ws.Cells.Item($r, 11).Value = $oldJ    # K (new): legislator_id
    $ws.Cells.Item($r, 10).Value = $oldI    # J: legislator_name
    $ws.Cells.Item($r, 9).Value = $oldH     # I: date (kept as text, see NumberFormat above)
    $ws.Cells.Item($r, 8).Value = "stock"   # H: property_category
}

# --- cleanup the two full-width-comma text values --------------------------
$ws.Range("G2").Value = "4797260"
$ws.Range("D17").Value = "2286"
